$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 835.5
$ws.Range("I101").Value = 314
$ws.Range("J101").Value = 1148.4
$ws.Range("K101").Value = 942
$ws.Range("L101").Value = 3445.2
$ws.Range("M101").Value = 680
$ws.Range("N101").Value = -6689.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3770.5833
$ws.Range("I113").Value = 2164.1428
$ws.Range("J113").Value = 6019.6
$ws.Range("K113").Value = 2164.1428
$ws.Range("L113").Value = 6019.6
$ws.Range("M113").Value = 1089.8572
$ws.Range("N113").Value = -12527.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7881.5
$ws.Range("I116").Value = 4312.1113
$ws.Range("K116").Value = 4312.1113
$ws.Range("M116").Value = -870.1112999999996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2657.162
$ws.Range("I132").Value = 2166.5356
$ws.Range("K132").Value = 6499.6068
$ws.Range("M132").Value = -3969.6068

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1074.4286
$ws.Range("I135").Value = 1059.0435
$ws.Range("J135").Value = 1310.3334
$ws.Range("K135").Value = 9531.3915
$ws.Range("L135").Value = 11793.0006
$ws.Range("M135").Value = -6996.3915
$ws.Range("N135").Value = -16863.0006

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1986.8422
$ws.Range("I138").Value = 1658.875
$ws.Range("K138").Value = 4976.625
$ws.Range("M138").Value = 163.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3055.6428
$ws.Range("I141").Value = 1335.1923
$ws.Range("J141").Value = 5851.375
$ws.Range("K141").Value = 4005.5769
$ws.Range("L141").Value = 17554.125
$ws.Range("M141").Value = 1174.4231
$ws.Range("N141").Value = -27914.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9486.905000000001
$ws.Range("I32").Value = 5627.8096
$ws.Range("J32").Value = 31589
$ws.Range("K32").Value = 5627.8096
$ws.Range("L32").Value = 31589
$ws.Range("M32").Value = -5340.8096
$ws.Range("N32").Value = -32163

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3645.1428
$ws.Range("I63").Value = 2902.5
$ws.Range("J63").Value = 4635.3335
$ws.Range("K63").Value = 2902.5
$ws.Range("L63").Value = 4635.3335
$ws.Range("M63").Value = -2216.5
$ws.Range("N63").Value = -6007.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3645.1428
$ws.Range("I66").Value = 2902.5
$ws.Range("J66").Value = 4635.3335
$ws.Range("K66").Value = 14512.5
$ws.Range("L66").Value = 23176.6675
$ws.Range("M66").Value = -11080.5
$ws.Range("N66").Value = -30040.6675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3614.5806
$ws.Range("I74").Value = 2294.7144
$ws.Range("K74").Value = 2294.7144
$ws.Range("M74").Value = -1420.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3614.5806
$ws.Range("I77").Value = 2294.7144
$ws.Range("K77").Value = 11473.572
$ws.Range("M77").Value = -7105.572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 963.4074000000001
$ws.Range("I97").Value = 829.5
$ws.Range("J97").Value = 1231.2222
$ws.Range("K97").Value = 829.5
$ws.Range("L97").Value = 1231.2222
$ws.Range("M97").Value = -333.5
$ws.Range("N97").Value = -2223.2222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 10099.857
$ws.Range("I132").Value = 9733
$ws.Range("J132").Value = 10375
$ws.Range("K132").Value = 29199
$ws.Range("L132").Value = 31125
$ws.Range("M132").Value = -26669
$ws.Range("N132").Value = -36185

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 102832.5
$ws.Range("J139").Value = 102832.5
$ws.Range("L139").Value = 102832.5
$ws.Range("N139").Value = -113112.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 11434.333
$ws.Range("J81").Value = 13379.4
$ws.Range("L81").Value = 13379.4
$ws.Range("N81").Value = -15501.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 11434.333
$ws.Range("J84").Value = 13379.4
$ws.Range("L84").Value = 40138.2
$ws.Range("N84").Value = -50746.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 22702.5
$ws.Range("I86").Value = 2398
$ws.Range("J86").Value = 43007
$ws.Range("K86").Value = 2398
$ws.Range("L86").Value = 43007
$ws.Range("M86").Value = -1275
$ws.Range("N86").Value = -45253

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 22702.5
$ws.Range("I89").Value = 2398
$ws.Range("J89").Value = 43007
$ws.Range("K89").Value = 11990
$ws.Range("L89").Value = 215035
$ws.Range("M89").Value = -6374
$ws.Range("N89").Value = -226267

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6434.311
$ws.Range("I31").Value = 2701.2917
$ws.Range("J31").Value = 10700.619
$ws.Range("K31").Value = 2701.2917
$ws.Range("L31").Value = 10700.619
$ws.Range("M31").Value = -2406.2917
$ws.Range("N31").Value = -11290.619

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6434.311
$ws.Range("I34").Value = 2701.2917
$ws.Range("J34").Value = 10700.619
$ws.Range("K34").Value = 2701.2917
$ws.Range("L34").Value = 10700.619
$ws.Range("M34").Value = -2499.2917
$ws.Range("N34").Value = -11104.619

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1724.4286
$ws.Range("I107").Value = 1444.64
$ws.Range("K107").Value = 1444.64
$ws.Range("M107").Value = 475.3599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 117773
$ws.Range("J138").Value = 117773
$ws.Range("L138").Value = 117773
$ws.Range("N138").Value = -128053

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 254.66667
$ws.Range("I13").Value = 254.66667
$ws.Range("K13").Value = 764.00001
$ws.Range("M13").Value = -596.00001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 151.25
$ws.Range("I15").Value = 101.42857
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 304.28571
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = -164.28571
$ws.Range("N15").Value = -1780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 307
$ws.Range("J36").Value = 239
$ws.Range("L36").Value = 717
$ws.Range("N36").Value = -1055

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 290.2857
$ws.Range("I40").Value = 302.6
$ws.Range("J40").Value = 283.44446
$ws.Range("K40").Value = 1210.4
$ws.Range("L40").Value = 1133.77784
$ws.Range("M40").Value = -1141.4
$ws.Range("N40").Value = -1271.77784

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1545.4286
$ws.Range("I117").Value = 1404.5
$ws.Range("J117").Value = 1733.3334
$ws.Range("K117").Value = 4213.5
$ws.Range("L117").Value = 5200.0002
$ws.Range("M117").Value = -771.5
$ws.Range("N117").Value = -12084.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 12823899
$ws.Range("I129").Value = 1854.8334
$ws.Range("J129").Value = 23814222
$ws.Range("K129").Value = 5564.5002
$ws.Range("L129").Value = 71442666
$ws.Range("M129").Value = -564.5002000000004
$ws.Range("N129").Value = -71452666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5846
$ws.Range("I131").Value = 3142.8333
$ws.Range("J131").Value = 8163
$ws.Range("K131").Value = 9428.499899999999
$ws.Range("L131").Value = 24489
$ws.Range("M131").Value = -4388.499899999999
$ws.Range("N131").Value = -34569

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3449.5
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 6386
$ws.Range("I17").Value = 10350
$ws.Range("J17").Value = 440
$ws.Range("K17").Value = 10350
$ws.Range("L17").Value = 440
$ws.Range("M17").Value = -10182
$ws.Range("N17").Value = -776

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 19275
$ws.Range("J58").Value = 26000
$ws.Range("L58").Value = 26000
$ws.Range("N58").Value = -26554

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3762.4
$ws.Range("I113").Value = 3699.6667
$ws.Range("J113").Value = 3856.5
$ws.Range("K113").Value = 3699.6667
$ws.Range("L113").Value = 3856.5
$ws.Range("M113").Value = -1529.6667
$ws.Range("N113").Value = -8196.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3180.8823
$ws.Range("I22").Value = 2307.9
$ws.Range("J22").Value = 4428
$ws.Range("K22").Value = 2307.9
$ws.Range("L22").Value = 4428
$ws.Range("M22").Value = -2012.9
$ws.Range("N22").Value = -5018

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3180.8823
$ws.Range("I27").Value = 2307.9
$ws.Range("J27").Value = 4428
$ws.Range("K27").Value = 2307.9
$ws.Range("L27").Value = 4428
$ws.Range("M27").Value = -2200.9
$ws.Range("N27").Value = -4642

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1295.5333
$ws.Range("I55").Value = 227.5
$ws.Range("J55").Value = 2007.5555
$ws.Range("K55").Value = 227.5
$ws.Range("L55").Value = 2007.5555
$ws.Range("M55").Value = -54.5
$ws.Range("N55").Value = -2353.5555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7152.8823
$ws.Range("I136").Value = 2825.25
$ws.Range("J136").Value = 9765.791999999999
$ws.Range("K136").Value = 8475.75
$ws.Range("L136").Value = 29297.376
$ws.Range("M136").Value = -5925.75
$ws.Range("N136").Value = -34397.376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1582
$ws.Range("I126").Value = 1582
$ws.Range("K126").Value = 4746
$ws.Range("M126").Value = -2276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1758.8451
$ws.Range("I136").Value = 1722.2449
$ws.Range("J136").Value = 1840.3636
$ws.Range("K136").Value = 5166.7347
$ws.Range("L136").Value = 5521.0908
$ws.Range("M136").Value = -2616.7347
$ws.Range("N136").Value = -10621.0908
